$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table by one more year column (R), inserting it so the new
# cells inherit column Q's formatting (border/number-format/font) exactly,
# the same way Excel does when a user inserts a column next to existing
# data.
$ws.Range("R2:R5").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftToRight) | Out-Null

# Fill in the new year's data (2021)
$ws.Range("R3").Value = 2021
$ws.Range("R4").Value = 202551
$ws.Range("R5").Value = 2.9794303052841493

# Update the active selection to the newly added cell
$ws.Range("R2").Select() | Out-Null
